$d = $word.ActiveDocument

$d.Content.Find.Execute("<id>p003r_a1</id>", $true, $false, $false, $false, $false, $true, 1, $false, "<id>p003r_1</id>", 2)
$d.Content.Find.Execute("<id>p003r_a2</id>", $true, $false, $false, $false, $false, $true, 1, $false, "<id>p003r_2</id>", 2)
$d.Content.Find.Execute("<id>p003r_a3</id>", $true, $false, $false, $false, $false, $true, 1, $false, "<id>p003r_3</id>", 2)
